# Implementation of document merging completed
#
# 1) Rename the two surviving pictures' internal <pic:cNvPr name="..."/>
#    (the "file name" shown in the drawing's picture-properties, which is
#    distinct from the wp:docPr "Picture N" display name and is not part
#    of the document's visible text, so Find/Replace cannot reach it).
# 2) Remove the trailing "Images:" + picture paragraph pairs that were
#    merged away, leaving only the first two image blocks.

$d = $word.ActiveDocument

function Rename-InlinePicture($doc, $oldName, $newName) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $r = $doc.Paragraphs.Item($i).Range
        $owx = $r.WordOpenXML
        $needle = 'name="' + $oldName + '"'
        if ($owx.Contains($needle)) {
            $newOwx = $owx.Replace($needle, 'name="' + $newName + '"')
            $r.InsertXML($newOwx)
            return $true
        }
    }
    return $false
}

# --- Step 1: rename the two pictures that remain in the final document ---
Rename-InlinePicture $d "Healthcare_AU_Implementation Approach (3).png" "Healthcare_AU_Implementation Approach (2).png" | Out-Null

Rename-InlinePicture $d "Healthcare_AU_Implementation Approach_Figure 6 KPMG Powered Enterprise Benefits to SJGHC (3).png" "Healthcare_AU_Implementation Approach_Figure 6 KPMG Powered Enterprise Benefits to SJGHC (1).png" | Out-Null

# --- Step 2: drop the extra "Images:" + picture paragraphs (Figure 7,
#     Figure 8, and the duplicate first picture) that trail the second
#     surviving image, right up to (but not including) the final
#     section-properties paragraph. ---
$anchorText = "Healthcare_AU_Implementation Approach_Figure 6 KPMG Powered Enterprise Benefits to SJGHC (1).png"
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $owx = $d.Paragraphs.Item($i).Range.WordOpenXML
    if ($owx.Contains($anchorText)) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0 -and $anchorIndex -lt $count) {
    $startPos = $d.Paragraphs.Item($anchorIndex + 1).Range.Start
    $endPos = $d.Paragraphs.Item($count).Range.End
    $d.Range($startPos, $endPos).Delete()
}
